$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: scoped, non-wrapping replace. Searches forward from the start of
# $scopeRange for $old and substitutes $new, preserving the formatting of the
# matched run. Returns $true/$false for whether a match was found.
# ---------------------------------------------------------------------------
function Replace-InRange($scopeRange, $old, $new) {
    $found = $scopeRange.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    return $found
}

# ---------------------------------------------------------------------------
# Helper used for paragraphs whose runs are being reshuffled (text moved from
# one run to another). Searches forward from $cursor (a Range) for the exact
# original run text, replaces just that run's text (preserving its rPr) and
# returns a new cursor range positioned right after the replaced run so the
# next call only ever looks at not-yet-processed (still-original) text.
# ---------------------------------------------------------------------------
function Replace-NextRun($cursor, $paraEnd, $oldText, $newText) {
    $found = $cursor.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output ("NOT FOUND: " + $oldText)
        return $cursor
    }
    $cursor.Text = $newText
    $newCursor = $d.Range($cursor.End, $paraEnd)
    return $newCursor
}

# 1) Language list line (hyperlink "英文" -> "英语", then the rest of the list)
Replace-InRange $d.Content "英文" "英语" | Out-Null
Replace-InRange $d.Content " / 葡萄牙文 / 法文 / 泰文 / 越南文 / 西班牙文" " / 葡萄牙语 / 法语 / 泰语 / 越南语 / 西班牙语" | Out-Null

# 2) Standalone "英文" heading paragraph (paragraph 3) - both occurrences of
#    "英文" already map to "英语" identically, handled by the global replace above.

# 3) Table cell contents
$tbl = $d.Tables.Item(1)
$cell = $tbl.Cell(1, 1)
Replace-InRange $cell.Range "簡介" "简要" | Out-Null
Replace-InRange $cell.Range "發送給在目標國家中已正確提交文檔的合作夥伴的驗證郵件。 將通過 customer.io 發送" "已发送给在目标国家提交正确文件的合作伙伴的验证电子邮件。 将通过 customer.io 发送" | Out-Null
Replace-InRange $cell.Range "目標受眾" "目标受众" | Out-Null
Replace-InRange $cell.Range "未按時提交文檔的邀請合作夥伴" "未按时提交文件的被邀请合作伙伴" | Out-Null

# 4) Subject line heading + placeholder
Replace-InRange $d.Content "主旨行" "主题行" | Out-Null
Replace-InRange $d.Content "[活動名稱]" "[活动名称]" | Out-Null

# 5) Big green header ("您的文檔已通過驗證！" -> "文件已验证!")
Replace-InRange $d.Content "您的文檔已通過驗證！" "文件已验证!" | Out-Null

# 6) Greeting paragraph: "[合作夥伴姓名]" placeholder + trailing comma
Replace-InRange $d.Content "[合作夥伴姓名]" "[合作伙伴姓名]" | Out-Null
$p16 = $d.Paragraphs.Item(16)
Replace-InRange $p16.Range ", " "， " | Out-Null

# 7) "我們已經審查了您傳送的 [活動名稱] 文檔，所有文檔均已通過驗證！ " paragraph
$p18 = $d.Paragraphs.Item(18)
Replace-InRange $p18.Range "我們已經審查了您傳送的 " "已经审查了您发送的 " | Out-Null
Replace-InRange $p18.Range " 文檔，所有文檔均已通過驗證！ " " 文件，所有文件均已通过验证！ " | Out-Null

# 8) "我們將很快傳送..." paragraph
Replace-InRange $d.Content "我們將很快傳送有關該活動的更多詳細資訊，包括議程和旅行安排，請務必定期查看收件箱。" "我们将很快发送有关该活动的更多详细信息，包括议程和旅行安排，请务必定期查看收件箱。" | Out-Null

# 9) Live-chat / WhatsApp contact paragraph (paragraph 20)
$p20 = $d.Paragraphs.Item(20)
Replace-InRange $p20.Range "如有任何疑問，請通過 " "If you have any questions, please contact us via " | Out-Null

$liveChatRange = $p20.Range.Duplicate()
$foundLc = $liveChatRange.Find.Execute("即時聊天", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundLc) {
    $liveChatRange.Text = "live chat"
    # The hyperlink run loses its direct character formatting when its text
    # is replaced; restore it explicitly (color 1155cc, single underline).
    $liveChatRange.Font.Color = 13391121
    $liveChatRange.Font.Underline = 1
}

Replace-InRange $p20.Range " 或 " " or " | Out-Null
Replace-InRange $p20.Range " 聯繫我們。 " ". " | Out-Null

# 10) Region-manager contact paragraph (paragraph 21) - runs are reshuffled,
#     so replace run-by-run from a forward-only, non-wrapping cursor.
$p21 = $d.Paragraphs.Item(21)
$paraEnd = $p21.Range.End
$cursor = $p21.Range.Duplicate()
$cursor = Replace-NextRun $cursor $paraEnd "如有任何疑問，請聯繫您的區域經理 " "如有任何疑问，请通过 "
$cursor = Replace-NextRun $cursor $paraEnd "[NAME]" "[电子邮件地址]"
$cursor = Replace-NextRun $cursor $paraEnd "，電子郵件地址為 " " 或 "
$cursor = Replace-NextRun $cursor $paraEnd "[電子郵件地址]" "[WHATSAPP 号码]"
$cursor = Replace-NextRun $cursor $paraEnd " 或 " " (WhatsApp) 联系您的区域经理 "
$cursor = Replace-NextRun $cursor $paraEnd "[WHATSAPP 號碼]" "[NAME]"
$cursor = Replace-NextRun $cursor $paraEnd " (WhatsApp)。 " "。 "

# 11) Comment text
$comment = $d.Comments.Item(1)
Replace-InRange $comment.Range "選擇任一項" "选择任一" | Out-Null

Write-Output "done"
